$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 122, shifting existing rows (122:195) down to (123:196).
$ws.Rows("122:122").Insert()

# Populate the newly inserted row 122 with the new price-report entry.
$ws.Range("A122").Value = 3
$ws.Range("B122").Value = "Femacal de La Calera"
$ws.Range("C122").Value = "Coquimbo"
$ws.Range("D122").Value = 44438
$ws.Range("E122").Value = 5
$ws.Range("F122").Value = 100112031
$ws.Range("G122").Value = "Poroto verde"
$ws.Range("H122").Value = "Magnum"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 85
$ws.Range("K122").Value = 27000
$ws.Range("L122").Value = 28000
$ws.Range("M122").Value = 27471
$ws.Range("N122").Value = "`$/malla 25 kilos"
$ws.Range("O122").Value = "Región de Arica y Parinacota"
$ws.Range("P122").Value = 1099
$ws.Range("Q122").Value = 25
$ws.Range("R122").Value = "Hortaliza"
